$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to existing data (rows 15 and 18) ---
$ws.Range("G15").Value = 15
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 17

# --- Register both numFmt codes (164 lowercase, 165 uppercase) while only
#     the uppercase one (165) ends up referenced by a cellXf, matching the
#     target style table. We do this by touching a scratch cell with the
#     lowercase code first, then the uppercase code (collapsing onto a
#     single new cellXf), then clearing the scratch cell entirely so it
#     leaves no trace in the used range / sheet data.
$scratch = $ws.Range("ZZ1")
$scratch.NumberFormat = "yyyy-mm-dd h:mm:ss"
$scratch.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$scratch.Clear()

# --- New column H: header ---
$ws.Range("A1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "dt_insertion"

# --- New column H: timestamp values for every data row ---
$dtValue = 45489.94295138889
for ($r = 2; $r -le 21; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $cell.Value = $dtValue
    $cell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

$excel.CutCopyMode = $false
